$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.141153333333333
$ws.Range("H2").Value = 12.42346
$ws.Range("I2").Value = 0.2530231305454066
$ws.Range("J2").Value = 0.2530231305454066
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.441151666666667
$ws.Range("N2").Value = 7.323455
$ws.Range("O2").Value = 0.1131710230270566
$ws.Range("P2").Value = 0.1131710230270566
$ws.Range("Q2").Value = 10.10918336158889
$ws.Range("R2").Value = 90.9826502543
$ws.Range("S2").Value = 0.02863488653333216
$ws.Range("T2").Value = 0.02863488653333216

$ws.Range("G3").Value = 4.141153333333333
$ws.Range("H3").Value = 12.42346
$ws.Range("I3").Value = 0.2530231305454066
$ws.Range("J3").Value = 0.2530231305454066
$ws.Range("O3").Value = 0.3310039188305578
$ws.Range("P3").Value = 0.3310039188305577
$ws.Range("Q3").Value = 29.56745657466222
$ws.Range("R3").Value = 266.10710917196
$ws.Range("S3").Value = 0.08375164776530541
$ws.Range("T3").Value = 0.0837516477653054

$ws.Range("G4").Value = 4.141153333333333
$ws.Range("H4").Value = 12.42346
$ws.Range("I4").Value = 0.2530231305454066
$ws.Range("J4").Value = 0.2530231305454066
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1123943333333333
$ws.Range("N4").Value = 0.337183
$ws.Range("O4").Value = 0.005210565922413945
$ws.Range("P4").Value = 0.005210565922413945
$ws.Range("Q4").Value = 0.465442168131111
$ws.Range("R4").Value = 4.18897951318
$ws.Range("S4").Value = 0.001318393701602391
$ws.Range("T4").Value = 0.001318393701602391

$ws.Range("G5").Value = 4.141153333333333
$ws.Range("H5").Value = 12.42346
$ws.Range("I5").Value = 0.2530231305454066
$ws.Range("J5").Value = 0.2530231305454066
$ws.Range("M5").Value = 11.877011
$ws.Range("N5").Value = 35.631033
$ws.Range("O5").Value = 0.5506144922199717
$ws.Range("P5").Value = 0.5506144922199717
$ws.Range("Q5").Value = 49.18452369268666
$ws.Range("R5").Value = 442.66071323418
$ws.Range("S5").Value = 0.1393182025451667
$ws.Range("T5").Value = 0.1393182025451667

$ws.Range("I6").Value = 0.3583796455306321
$ws.Range("J6").Value = 0.358379645530632
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.441151666666667
$ws.Range("N6").Value = 7.323455
$ws.Range("O6").Value = 0.1131710230270566
$ws.Range("P6").Value = 0.1131710230270566
$ws.Range("Q6").Value = 14.31855475790278
$ws.Range("R6").Value = 128.866992821125
$ws.Range("S6").Value = 0.04055819111677554
$ws.Range("T6").Value = 0.04055819111677553

$ws.Range("I7").Value = 0.3583796455306321
$ws.Range("J7").Value = 0.358379645530632
$ws.Range("O7").Value = 0.3310039188305578
$ws.Range("P7").Value = 0.3310039188305577
$ws.Range("S7").Value = 0.1186250670997454
$ws.Range("T7").Value = 0.1186250670997454

$ws.Range("I8").Value = 0.3583796455306321
$ws.Range("J8").Value = 0.358379645530632
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1123943333333333
$ws.Range("N8").Value = 0.337183
$ws.Range("O8").Value = 0.005210565922413945
$ws.Range("P8").Value = 0.005210565922413945
$ws.Range("Q8").Value = 0.6592480255472223
$ws.Range("R8").Value = 5.933232229925001
$ws.Range("S8").Value = 0.0018673607682887
$ws.Range("T8").Value = 0.0018673607682887

$ws.Range("I9").Value = 0.3583796455306321
$ws.Range("J9").Value = 0.358379645530632
$ws.Range("M9").Value = 11.877011
$ws.Range("N9").Value = 35.631033
$ws.Range("O9").Value = 0.5506144922199717
$ws.Range("P9").Value = 0.5506144922199717
$ws.Range("Q9").Value = 69.66450904540835
$ws.Range("R9").Value = 626.9805814086751
$ws.Range("S9").Value = 0.1973290265458224
$ws.Range("T9").Value = 0.1973290265458224

$ws.Range("G10").Value = 0.467591
$ws.Range("H10").Value = 1.402773
$ws.Range("I10").Value = 0.02856965900840602
$ws.Range("J10").Value = 0.02856965900840601
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.441151666666667
$ws.Range("N10").Value = 7.323455
$ws.Range("O10").Value = 0.1131710230270566
$ws.Range("P10").Value = 0.1131710230270566
$ws.Range("Q10").Value = 1.141460548968334
$ws.Range("R10").Value = 10.273144940715
$ws.Range("S10").Value = 0.003233257537515472
$ws.Range("T10").Value = 0.003233257537515471

$ws.Range("G11").Value = 0.467591
$ws.Range("H11").Value = 1.402773
$ws.Range("I11").Value = 0.02856965900840602
$ws.Range("J11").Value = 0.02856965900840601
$ws.Range("O11").Value = 0.3310039188305578
$ws.Range("P11").Value = 0.3310039188305577
$ws.Range("Q11").Value = 3.338557033355334
$ws.Range("R11").Value = 30.047013300198
$ws.Range("S11").Value = 0.00945666909143514
$ws.Range("T11").Value = 0.009456669091435136

$ws.Range("G12").Value = 0.467591
$ws.Range("H12").Value = 1.402773
$ws.Range("I12").Value = 0.02856965900840602
$ws.Range("J12").Value = 0.02856965900840601
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.1123943333333333
$ws.Range("N12").Value = 0.337183
$ws.Range("O12").Value = 0.005210565922413945
$ws.Range("P12").Value = 0.005210565922413945
$ws.Range("Q12").Value = 0.05255457871766667
$ws.Range("R12").Value = 0.4729912084590001
$ws.Range("S12").Value = 0.000148864091644187
$ws.Range("T12").Value = 0.0001488640916441869

$ws.Range("G13").Value = 0.467591
$ws.Range("H13").Value = 1.402773
$ws.Range("I13").Value = 0.02856965900840602
$ws.Range("J13").Value = 0.02856965900840601
$ws.Range("M13").Value = 11.877011
$ws.Range("N13").Value = 35.631033
$ws.Range("O13").Value = 0.5506144922199717
$ws.Range("P13").Value = 0.5506144922199717
$ws.Range("Q13").Value = 5.553583450501001
$ws.Range("R13").Value = 49.98225105450901
$ws.Range("S13").Value = 0.01573086828781122
$ws.Range("T13").Value = 0.01573086828781122

$ws.Range("G14").Value = 5.892462666666667
$ws.Range("H14").Value = 17.677388
$ws.Range("I14").Value = 0.3600275649155554
$ws.Range("J14").Value = 0.3600275649155554
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.441151666666667
$ws.Range("N14").Value = 7.323455
$ws.Range("O14").Value = 0.1131710230270566
$ws.Range("P14").Value = 0.1131710230270566
$ws.Range("Q14").Value = 14.38439505950445
$ws.Range("R14").Value = 129.45955553554
$ws.Range("S14").Value = 0.04074468783943343
$ws.Range("T14").Value = 0.04074468783943342

$ws.Range("G15").Value = 5.892462666666667
$ws.Range("H15").Value = 17.677388
$ws.Range("I15").Value = 0.3600275649155554
$ws.Range("J15").Value = 0.3600275649155554
$ws.Range("O15").Value = 0.3310039188305578
$ws.Range("P15").Value = 0.3310039188305577
$ws.Range("Q15").Value = 42.07164526174311
$ws.Range("R15").Value = 378.644807355688
$ws.Range("S15").Value = 0.1191705348740719
$ws.Range("T15").Value = 0.1191705348740719

$ws.Range("G16").Value = 5.892462666666667
$ws.Range("H16").Value = 17.677388
$ws.Range("I16").Value = 0.3600275649155554
$ws.Range("J16").Value = 0.3600275649155554
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.1123943333333333
$ws.Range("N16").Value = 0.337183
$ws.Range("O16").Value = 0.005210565922413945
$ws.Range("P16").Value = 0.005210565922413945
$ws.Range("Q16").Value = 0.6622794131115556
$ws.Range("R16").Value = 5.960514718004
$ws.Range("S16").Value = 0.001875947360878668
$ws.Range("T16").Value = 0.001875947360878667

$ws.Range("G17").Value = 5.892462666666667
$ws.Range("H17").Value = 17.677388
$ws.Range("I17").Value = 0.3600275649155554
$ws.Range("J17").Value = 0.3600275649155554
$ws.Range("M17").Value = 11.877011
$ws.Range("N17").Value = 35.631033
$ws.Range("O17").Value = 0.5506144922199717
$ws.Range("P17").Value = 0.5506144922199717
$ws.Range("Q17").Value = 69.98484390908935
$ws.Range("R17").Value = 629.863595181804
$ws.Range("S17").Value = 0.1982363948411714
$ws.Range("T17").Value = 0.1982363948411714

